# Update readme and template
# Swap the "Description"/"Date" column headers in both the Expenses (D2:E2)
# and Income (L2:M2) header rows, bold-format the now-active cell D5, and
# move the active selection to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Expenses header: swap D2 ("Description") and E2 ("Date") ---
$d2 = $ws.Range("D2").Value2
$e2 = $ws.Range("E2").Value2
$ws.Range("D2").Value = $e2
$ws.Range("E2").Value = $d2

# --- Income header: swap L2 ("Description") and M2 ("Date") ---
$l2 = $ws.Range("L2").Value2
$m2 = $ws.Range("M2").Value2
$ws.Range("L2").Value = $m2
$ws.Range("M2").Value = $l2

# --- Bold-format D5 (empty cell, matches the header style) ---
$ws.Range("D5").Font.Bold = $true

# --- Move the active selection to D5 ---
$ws.Range("D5").Select() | Out-Null
